$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.834.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "'2.495.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'532.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "'134.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("D11").Value = "'5.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").Value = "'2.937.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").Value = "'58.765.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "'22.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.24%  "
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").Value = "'2.508.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "'11.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").Value = "'322.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "'5.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.46%  "
$ws.Range("D23").Value = "'64.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.53%  "
$ws.Range("D24").Value = "'0.418"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("D29").Value = "'169.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.83%  "
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("D31").Value = "'6.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.11%  "
$ws.Range("E32").Value = "  +0.65%  "
$ws.Range("D34").Value = "'18.28"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("E37").Value = "  -2.63%  "
$ws.Range("E38").Value = "  -1.09%  "
$ws.Range("D39").Value = "'0.797"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.14%  "
$ws.Range("D40").Value = "'280.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.16%  "
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("E42").Value = "  -4.71%  "
$ws.Range("D44").Value = "'129.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.45%  "
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").Value = "'0.0924"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("E47").Value = "  -2.20%  "
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("D49").Value = "'17.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("D50").Value = "'1.749.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("E51").Value = "  -0.43%  "
